# Weekly crime-data refresh: advance report one week (Vol 30 No 42 -> No 43;
# week of 10/16-10/22/2023 -> 10/23-10/29/2023) and update all Week-to-Date,
# 28-Day, Year-to-Date counts and percentage changes for rows 15-30.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and week-covering dates) ---
$ws.Cells.Item(8,1).Value = "Volume 30   Number  43"
$ws.Cells.Item(9,3).Value = "Report Covering the Week  10/23/2023  Through  10/29/2023"

# --- Cells whose type flips between a numeric value and the ")" / "***.*" text
# placeholders: copy format+value from a same-shaped donor cell elsewhere in the
# table (fixes the style index + shared-string type), then set the real value. ---
$ws.Cells.Item(20,4).Copy($ws.Cells.Item(15,4))
$ws.Cells.Item(20,5).Copy($ws.Cells.Item(15,5))
$ws.Cells.Item(22,6).Copy($ws.Cells.Item(15,6))
$ws.Cells.Item(15,3).Copy($ws.Cells.Item(16,3))
$ws.Cells.Item(17,4).Copy($ws.Cells.Item(16,4))
$ws.Cells.Item(16,4).Value = 2
$ws.Cells.Item(17,5).Copy($ws.Cells.Item(16,5))
$ws.Cells.Item(16,5).Value = -100
$ws.Cells.Item(17,3).Copy($ws.Cells.Item(26,3))
$ws.Cells.Item(26,3).Value = 1
$ws.Cells.Item(20,4).Copy($ws.Cells.Item(26,4))
$ws.Cells.Item(20,5).Copy($ws.Cells.Item(26,5))
$ws.Cells.Item(17,3).Copy($ws.Cells.Item(28,3))
$ws.Cells.Item(28,3).Value = 1
$ws.Cells.Item(17,3).Copy($ws.Cells.Item(29,3))
$ws.Cells.Item(29,3).Value = 1
$ws.Cells.Item(17,4).Copy($ws.Cells.Item(30,4))
$ws.Cells.Item(30,4).Value = 1
$ws.Cells.Item(17,5).Copy($ws.Cells.Item(30,5))
$ws.Cells.Item(30,5).Value = -100

# --- Remaining cells: same type before/after, just update the numbers. ---
$ws.Cells.Item(15,7).Value = 1
$ws.Cells.Item(15,8).Value = -100
$ws.Cells.Item(15,12).Value = -27.777777777777
$ws.Cells.Item(15,14).Value = -63.888888888888
$ws.Cells.Item(16,6).Value = 6
$ws.Cells.Item(16,8).Value = 20
$ws.Cells.Item(16,9).Value = 84
$ws.Cells.Item(16,10).Value = 77
$ws.Cells.Item(16,11).Value = 9.090909090909
$ws.Cells.Item(16,12).Value = -1.176470588235
$ws.Cells.Item(16,13).Value = -19.230769230769
$ws.Cells.Item(16,14).Value = -85.237258347978
$ws.Cells.Item(17,3).Value = 8
$ws.Cells.Item(17,4).Value = 10
$ws.Cells.Item(17,5).Value = -20
$ws.Cells.Item(17,6).Value = 26
$ws.Cells.Item(17,7).Value = 19
$ws.Cells.Item(17,8).Value = 36.842105263157
$ws.Cells.Item(17,9).Value = 227
$ws.Cells.Item(17,10).Value = 200
$ws.Cells.Item(17,11).Value = 13.5
$ws.Cells.Item(17,12).Value = -8.467741935483
$ws.Cells.Item(17,13).Value = 69.402985074626
$ws.Cells.Item(17,14).Value = -46.462264150943
$ws.Cells.Item(18,3).Value = 2
$ws.Cells.Item(18,4).Value = 4
$ws.Cells.Item(18,5).Value = -50
$ws.Cells.Item(18,6).Value = 7
$ws.Cells.Item(18,7).Value = 8
$ws.Cells.Item(18,8).Value = -12.5
$ws.Cells.Item(18,9).Value = 61
$ws.Cells.Item(18,10).Value = 65
$ws.Cells.Item(18,11).Value = -6.153846153846
$ws.Cells.Item(18,12).Value = 3.389830508474
$ws.Cells.Item(18,13).Value = -46.956521739130
$ws.Cells.Item(18,14).Value = -90.895522388059
$ws.Cells.Item(19,3).Value = 3
$ws.Cells.Item(19,4).Value = 1
$ws.Cells.Item(19,5).Value = 200
$ws.Cells.Item(19,7).Value = 9
$ws.Cells.Item(19,8).Value = 88.888888888888
$ws.Cells.Item(19,9).Value = 159
$ws.Cells.Item(19,10).Value = 145
$ws.Cells.Item(19,11).Value = 9.655172413793
$ws.Cells.Item(19,12).Value = 8.163265306122
$ws.Cells.Item(19,13).Value = 84.883720930232
$ws.Cells.Item(19,14).Value = -40.892193308550
$ws.Cells.Item(20,3).Value = 2
$ws.Cells.Item(20,7).Value = 3
$ws.Cells.Item(20,8).Value = 66.666666666666
$ws.Cells.Item(20,9).Value = 56
$ws.Cells.Item(20,11).Value = 9.803921568627
$ws.Cells.Item(20,12).Value = 51.351351351351
$ws.Cells.Item(20,13).Value = -18.840579710144
$ws.Cells.Item(20,14).Value = -85.750636132315
$ws.Cells.Item(21,3).Value = 15
$ws.Cells.Item(21,4).Value = 17
$ws.Cells.Item(21,5).Value = -11.764705882352
$ws.Cells.Item(21,6).Value = 61
$ws.Cells.Item(21,7).Value = 46
$ws.Cells.Item(21,8).Value = 32.608695652173
$ws.Cells.Item(21,9).Value = 600
$ws.Cells.Item(21,10).Value = 562
$ws.Cells.Item(21,11).Value = 6.761565836298
$ws.Cells.Item(21,12).Value = -0.826446280991
$ws.Cells.Item(21,13).Value = 14.068441064638
$ws.Cells.Item(21,14).Value = -74.747474747474
$ws.Cells.Item(22,7).Value = 1
$ws.Cells.Item(23,3).Value = 1
$ws.Cells.Item(23,4).Value = 4
$ws.Cells.Item(23,5).Value = -75
$ws.Cells.Item(23,7).Value = 10
$ws.Cells.Item(23,8).Value = -40
$ws.Cells.Item(23,9).Value = 71
$ws.Cells.Item(23,10).Value = 75
$ws.Cells.Item(23,11).Value = -5.333333333333
$ws.Cells.Item(23,12).Value = -13.414634146341
$ws.Cells.Item(23,13).Value = 108.823529411765
$ws.Cells.Item(24,3).Value = 15
$ws.Cells.Item(24,4).Value = 10
$ws.Cells.Item(24,5).Value = 50
$ws.Cells.Item(24,6).Value = 49
$ws.Cells.Item(24,7).Value = 40
$ws.Cells.Item(24,8).Value = 22.5
$ws.Cells.Item(24,9).Value = 484
$ws.Cells.Item(24,10).Value = 489
$ws.Cells.Item(24,11).Value = -1.022494887525
$ws.Cells.Item(24,12).Value = -5.836575875486
$ws.Cells.Item(24,13).Value = 55.128205128205
$ws.Cells.Item(25,3).Value = 9
$ws.Cells.Item(25,4).Value = 5
$ws.Cells.Item(25,5).Value = 80
$ws.Cells.Item(25,6).Value = 31
$ws.Cells.Item(25,7).Value = 24
$ws.Cells.Item(25,8).Value = 29.166666666666
$ws.Cells.Item(25,9).Value = 357
$ws.Cells.Item(25,10).Value = 329
$ws.Cells.Item(25,11).Value = 8.510638297872
$ws.Cells.Item(25,12).Value = -3.773584905660
$ws.Cells.Item(25,13).Value = 3.478260869565
$ws.Cells.Item(26,7).Value = 2
$ws.Cells.Item(26,8).Value = -50
$ws.Cells.Item(26,9).Value = 26
$ws.Cells.Item(26,11).Value = -10.344827586206
$ws.Cells.Item(26,12).Value = 13.043478260869
$ws.Cells.Item(27,4).Value = 1
$ws.Cells.Item(27,7).Value = 7
$ws.Cells.Item(27,8).Value = -85.714285714285
$ws.Cells.Item(27,10).Value = 45
$ws.Cells.Item(27,11).Value = -24.444444444444
$ws.Cells.Item(27,12).Value = -19.047619047619
$ws.Cells.Item(28,9).Value = 12
$ws.Cells.Item(28,11).Value = -52
$ws.Cells.Item(28,12).Value = -53.846153846153
$ws.Cells.Item(28,13).Value = -42.857142857142
$ws.Cells.Item(28,14).Value = -68.421052631578
$ws.Cells.Item(29,9).Value = 10
$ws.Cells.Item(29,11).Value = -47.368421052631
$ws.Cells.Item(29,12).Value = -56.521739130434
$ws.Cells.Item(29,13).Value = -47.368421052631
$ws.Cells.Item(29,14).Value = -72.222222222222
$ws.Cells.Item(30,10).Value = 4
$ws.Cells.Item(30,11).Value = -50
